# Update "Forecast Comparison" sheet with corrected forecast output:
#  - insert a new "Week_Start_Date" column right after "Week" (new col B)
#  - shift ASIN / MyForecast / Amazon Mean / P70 / P80 / P90 / Product Title /
#    is_holiday_week one column to the right (C..J)
#  - strip the leading zero from single-digit week labels (W01 -> W1, ... W09 -> W9)
#  - re-type the is_holiday_week column as a boolean instead of a number

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Week"
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"
$ws.Cells.Item(1, 3).Value = "ASIN"
$ws.Cells.Item(1, 4).Value = "MyForecast"
$ws.Cells.Item(1, 5).Value = "Amazon Mean Forecast"
$ws.Cells.Item(1, 6).Value = "Amazon P70 Forecast"
$ws.Cells.Item(1, 7).Value = "Amazon P80 Forecast"
$ws.Cells.Item(1, 8).Value = "Amazon P90 Forecast"
$ws.Cells.Item(1, 9).Value = "Product Title"
$ws.Cells.Item(1, 10).Value = "is_holiday_week"

# --- Data rows -----------------------------------------------------------
$weeks      = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$startDates = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")
$asin       = "B07HS59X7P"
$myForecast = @(2,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0)
$meanFc     = @(2,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0)
$p70        = @(1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0)
$p80        = @(3,2,2,1,0,0,0,0,0,0,0,0,0,0,0,0)
$p90        = @(5,3,3,1,0,0,0,0,0,0,0,0,0,0,0,0)
$title      = "Z390 UD"
$holiday    = @($false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false)

for ($i = 0; $i -lt $weeks.Length; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 1).Value = $weeks[$i]

    # Force the start date to be stored as plain text, not auto-converted
    # to a date serial number.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $startDates[$i]

    $ws.Cells.Item($row, 3).Value = $asin
    $ws.Cells.Item($row, 4).Value = $myForecast[$i]
    $ws.Cells.Item($row, 5).Value = $meanFc[$i]
    $ws.Cells.Item($row, 6).Value = $p70[$i]
    $ws.Cells.Item($row, 7).Value = $p80[$i]
    $ws.Cells.Item($row, 8).Value = $p90[$i]
    $ws.Cells.Item($row, 9).Value = $title
    $ws.Cells.Item($row, 10).Value = $holiday[$i]
}

Write-Output "Forecast Comparison sheet updated"
